$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "moduleCacheConfig"

# Column widths (A/B/C)
$ws.Columns.Item(1).ColumnWidth = 34.44140625
$ws.Columns.Item(2).ColumnWidth = 27.88671875
$ws.Columns.Item(3).ColumnWidth = 26.6640625

# Header row (row 1) values
$ws.Range("A1").Value = "test-id"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "moduleName"
$ws.Range("D1").Value = "rspStatus"
$ws.Range("E1").Value = "rspCode"
$ws.Range("F1").Value = "rspMessage"

# Data row (row 2) values
$ws.Range("A2").Value = "jinzu-connector-configure-test-1"
$ws.Range("B2").Value = "get module cache config"
$ws.Range("C2").Value = "data-layer-api-engine"
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "Operate success."

# Header row styling: 10pt font, shaded theme fill, thin border all around,
# vertically centered text
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Name = "等线"
$headerRange.Font.Size = 10
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 4
$headerRange.Borders.LineStyle = 1
$headerRange.VerticalAlignment = -4108

# Selection / active cell, matching the saved sheet view
$ws.Range("C8").Select() | Out-Null

# Page orientation
$ws.PageSetup.Orientation = 1 | Out-Null
